$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear stale data (rows 2-8) before rewriting the refreshed dataset (now only 6 data rows, 2-7)
$ws.Range("A2:AQ8").ClearContents()

# Row 2
$ws.Range("A2").Value = 'Luxembourg'
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = '5'
$ws.Range("C2").Value = 'Investments & Asset Management'
$ws.Range("D2").Value = 0.03351999999999999
$ws.Range("E2").Value = -0.0959
$ws.Range("G2").Value = -3.923281428028865
$ws.Range("H2").Value = -3.923281428028865
$ws.Range("I2").Value = 1.035067730092417
$ws.Range("J2").Value = 1.021252121754655
$ws.Range("K2").Value = 61.62
$ws.Range("L2").Value = -0.7800987466767947
$ws.Range("M2").Value = 229.117
$ws.Range("N2").Value = 0.02462035246077799
$ws.Range("O2").Value = 3.718224602401818
$ws.Range("P2").Value = 184.5
$ws.Range("Q2").Value = 0.01982591876208897
$ws.Range("R2").Value = 2.994157740993184
$ws.Range("S2").Value = 44.617
$ws.Range("T2").Value = 0.1947345679281764
$ws.Range("U2").Value = 122.46
$ws.Range("V2").Value = 0.01315925209542231
$ws.Range("W2").Value = 0.00969353964717924
$ws.Range("X2").Value = 0.03525885212668903
$ws.Range("Y2").Value = -0.02556531247950979
$ws.Range("Z2").Value = -0.007505458774611996
$ws.Range("AA2").Value = 0.009443189854483775
$ws.Range("AB2").Value = 0.03524549617030377
$ws.Range("AC2").Value = -0.02578601808003933
$ws.Range("AD2").Value = 25.17
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 25.17
$ws.Range("AG2").Value = -97.29000000000001
$ws.Range("AH2").Value = 0.002697410935606146
$ws.Range("AI2").Value = 0.002305711513992289
$ws.Range("AJ2").Value = -0.01056499770326137
$ws.Range("AK2").Value = -0.009013415898409382
$ws.Range("AL2").Value = 2.607
$ws.Range("AM2").Value = 1.295
$ws.Range("AN2").Value = 0.318204804045512
$ws.Range("AO2").Value = -31.36171845032604
$ws.Range("AP2").Value = -1.229962073324905
$ws.Range("AQ2").Value = -63.13513513513514

# Row 3
$ws.Range("A3").Value = 'Luxembourg'
$ws.Range("B3").Value = 'BBGI Global Infrastructure S.A. (LSE:BBGI)'
$ws.Range("C3").Value = 'Investments & Asset Management'
$ws.Range("D3").Value = 0.06519999999999999
$ws.Range("E3").Value = 0.0358
$ws.Range("G3").Value = 0.9697732997481108
$ws.Range("H3").Value = 0.9697732997481108
$ws.Range("I3").Value = 0.8035264483627204
$ws.Range("J3").Value = 0.7510987431288324
$ws.Range("K3").Value = 54.8
$ws.Range("L3").Value = 0.690176322418136
$ws.Range("M3").Value = 53.3
$ws.Range("N3").Value = 0.03375340383762903
$ws.Range("O3").Value = 0.9726277372262774
$ws.Range("P3").Value = 53.3
$ws.Range("Q3").Value = 0.03375340383762903
$ws.Range("R3").Value = 0.9726277372262774
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 18.9
$ws.Range("V3").Value = 0.01196884301184219
$ws.Range("W3").Value = 0.05036764705882352
$ws.Range("X3").Value = 0.03534687160382731
$ws.Range("Y3").Value = 0.01502077545499621
$ws.Range("Z3").Value = 0.07446309668948702
$ws.Range("AA3").Value = 0.05592913833295442
$ws.Range("AB3").Value = 0.03523845704962056
$ws.Range("AC3").Value = 0.02069068128333387
$ws.Range("AD3").Value = 8.220000000000001
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 8.220000000000001
$ws.Range("AG3").Value = -10.68
$ws.Range("AH3").Value = 0.005178539928936825
$ws.Range("AI3").Value = 0.00767994618431871
$ws.Range("AJ3").Value = -0.00680940054322184
$ws.Range("AK3").Value = -0.01015769150291986
$ws.Range("AL3").Value = 2.34
$ws.Range("AM3").Value = 2.258
$ws.Range("AN3").Value = 0.1288401253918496
$ws.Range("AO3").Value = 27.26495726495726
$ws.Range("AP3").Value = -0.167398119122257
$ws.Range("AQ3").Value = 28.25509300265722

# Row 4
$ws.Range("A4").Value = 'Luxembourg'
$ws.Range("B4").Value = 'Brederode SA (ENXTBR:BREB)'
$ws.Range("C4").Value = 'Investments & Asset Management'
$ws.Range("E4").Value = -0.0959
$ws.Range("K4").Value = 185.1
$ws.Range("M4").Value = 33.6
$ws.Range("N4").Value = 0.0120378331900258
$ws.Range("O4").Value = 0.1815235008103728
$ws.Range("P4").Value = 33.6
$ws.Range("Q4").Value = 0.0120378331900258
$ws.Range("R4").Value = 0.1815235008103728
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0.09
$ws.Range("V4").Value = 0.00003224419604471195
$ws.Range("W4").Value = 0.07228492209161558
$ws.Range("X4").Value = 0.03524549617030377
$ws.Range("Y4").Value = 0.03703942592131181
$ws.Range("Z4").Value = 0
$ws.Range("AA4").Value = 0.01308767570692983
$ws.Range("AB4").Value = 0.03524549617030377
$ws.Range("AC4").Value = -0.02215782046337394
$ws.Range("AD4").Value = 0
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 0
$ws.Range("AG4").Value = -0.09
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = -0.00003224523576641551
$ws.Range("AK4").Value = -0.00003353441562554727
$ws.Range("AL4").Value = 0
$ws.Range("AM4").Value = 0

# Row 5
$ws.Range("A5").Value = 'Luxembourg'
$ws.Range("B5").Value = 'Luxempart S.A. (BDL:LXMPR)'
$ws.Range("C5").Value = 'Investments & Asset Management'
$ws.Range("D5").Value = 0.00184
$ws.Range("E5").Value = -0.207
$ws.Range("G5").Value = 8.086805555555555
$ws.Range("H5").Value = 8.086805555555555
$ws.Range("I5").Value = 0.53125
$ws.Range("J5").Value = 0.5304580745341615
$ws.Range("K5").Value = 16.1
$ws.Range("L5").Value = 0.5590277777777778
$ws.Range("M5").Value = 33.517
$ws.Range("N5").Value = 0.02784266489450075
$ws.Range("O5").Value = 2.081801242236025
$ws.Range("P5").Value = 33.5
$ws.Range("Q5").Value = 0.02782854294733345
$ws.Range("R5").Value = 2.080745341614906
$ws.Range("S5").Value = 0.01700000000000301
$ws.Range("T5").Value = 0.000507205298803682
$ws.Range("U5").Value = 46.6
$ws.Range("V5").Value = 0.03871074929390264
$ws.Range("W5").Value = 0.00969353964717924
$ws.Range("X5").Value = 0.03548169027697537
$ws.Range("Y5").Value = -0.02578815062979613
$ws.Range("Z5").Value = 0.01780195326987267
$ws.Range("AA5").Value = 0.009443189854483775
$ws.Range("AB5").Value = 0.0352292079345231
$ws.Range("AC5").Value = -0.02578601808003933
$ws.Range("AD5").Value = 14.6
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 14.6
$ws.Range("AG5").Value = -32
$ws.Range("AH5").Value = 0.01198292843072883
$ws.Range("AI5").Value = 0.008901896225839888
$ws.Range("AJ5").Value = -0.0273084144051886
$ws.Range("AK5").Value = -0.02008158142453718
$ws.Range("AL5").Value = 0.255
$ws.Range("AM5").Value = -0.975
$ws.Range("AN5").Value = 0.9542483660130718
$ws.Range("AO5").Value = 60
$ws.Range("AP5").Value = -2.091503267973856
$ws.Range("AQ5").Value = -15.69230769230769

# Row 6
$ws.Range("A6").Value = 'Luxembourg'
$ws.Range("B6").Value = 'Reinet Investments S.C.A. (BDL:REINI)'
$ws.Range("C6").Value = 'Investments & Asset Management'
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 1.012472885032538
$ws.Range("J6").Value = 1.012472885032538
$ws.Range("K6").Value = -186.7
$ws.Range("L6").Value = 1.012472885032538
$ws.Range("M6").Value = 85.7
$ws.Range("N6").Value = 0.02501021420650208
$ws.Range("O6").Value = -0.4590251740760579
$ws.Range("P6").Value = 41.1
$ws.Range("Q6").Value = 0.01199439677814744
$ws.Range("R6").Value = -0.2201392608462775
$ws.Range("S6").Value = 44.6
$ws.Range("T6").Value = 0.5204200700116686
$ws.Range("U6").Value = 1.17
$ws.Range("V6").Value = 0.0003414463316406934
$ws.Range("W6").Value = -0.03608496492007963
$ws.Range("X6").Value = 0.03525885212668903
$ws.Range("Y6").Value = -0.07134381704676866
$ws.Range("Z6").Value = -0.03570059513898838
$ws.Range("AA6").Value = -0.03614588455775017
$ws.Range("AB6").Value = 0.0352479597472441
$ws.Range("AC6").Value = -0.07139384430499426
$ws.Range("AD6").Value = 2.35
$ws.Range("AE6").Value = 0
$ws.Range("AF6").Value = 2.35
$ws.Range("AG6").Value = 1.18
$ws.Range("AH6").Value = 0.0006853409935986235
$ws.Range("AI6").Value = 0.000439543996483648
$ws.Range("AJ6").Value = 0.000344246130148376
$ws.Range("AK6").Value = 0.0002207555076628352
$ws.Range("AL6").Value = 0
$ws.Range("AM6").Value = 0

# Row 7
$ws.Range("A7").Value = 'Luxembourg'
$ws.Range("B7").Value = 'NB Aurora S.A. SICAF-RAIF (BIT:NBA)'
$ws.Range("C7").Value = 'Investments & Asset Management'
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 2.745519713261649
$ws.Range("J7").Value = 2.745519713261649
$ws.Range("K7").Value = -7.68
$ws.Range("L7").Value = 2.752688172043011
$ws.Range("M7").Value = 23
$ws.Range("N7").Value = 0.07533573534228627
$ws.Range("O7").Value = -2.994791666666667
$ws.Range("P7").Value = 23
$ws.Range("Q7").Value = 0.07533573534228627
$ws.Range("R7").Value = -2.994791666666667
$ws.Range("S7").Value = 0
$ws.Range("T7").Value = 0
$ws.Range("U7").Value = 55.7
$ws.Range("V7").Value = 0.1824434981984933
$ws.Range("W7").Value = -0.03711938134364427
$ws.Range("X7").Value = 0.03524549617030377
$ws.Range("Y7").Value = -0.07236487751394804
$ws.Range("Z7").Value = -0.02417677642980936
$ws.Range("AA7").Value = -0.06637781629116117
$ws.Range("AB7").Value = 0.03524549617030377
$ws.Range("AC7").Value = -0.1016233124614649
$ws.Range("AD7").Value = 0
$ws.Range("AE7").Value = 0
$ws.Range("AF7").Value = 0
$ws.Range("AG7").Value = -55.7
$ws.Range("AH7").Value = 0
$ws.Range("AI7").Value = 0
$ws.Range("AJ7").Value = -0.2231570512820513
$ws.Range("AK7").Value = -0.4645537948290243
$ws.Range("AL7").Value = 0.012
$ws.Range("AM7").Value = 0.012
$ws.Range("AO7").Value = -638.3333333333334
$ws.Range("AQ7").Value = -638.3333333333334

